$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 with new iteration values
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "2.5"
$ws.Range("C2").Value = "11.125"
$ws.Range("D2").Value = "1.33256010122311"

$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "1.16743989877689"
$ws.Range("C3").Value = "-1.57631747829391"
$ws.Range("D3").Value = "0.236559671274285"

$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "1.40399957005117"
$ws.Range("C4").Value = "-0.636414848616419"
$ws.Range("D4").Value = "0.105985767951837"

$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "1.50998533800301"
$ws.Range("C5").Value = "-0.0671346294871952"
$ws.Range("D5").Value = "0.011293924179594"

$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "1.5212792621826"
$ws.Range("C6").Value = "-0.0005969755555189"
$ws.Range("D6").Value = "0.0001004368741777"

# Add new row 7 with the additional iteration
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "1.52137969905678"
$ws.Range("C7").Value = "-4.60512058374181e-08"
$ws.Range("D7").Value = "7.747786634482171e-09"
